# Sync attendance_reports: fix "Recorded By" column order
# Change "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in column G (Recorded By) of the active worksheet, for every matching cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
